# Code to write cell data in excel
# Populates the "Custom" worksheet (sheet2) with a Field/Override/Status
# table header plus three data rows (List One / Two / Three).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Custom")

# Header row
$ws.Range("A1").Value = "Field"
$ws.Range("B1").Value = "Override"
$ws.Range("C1").Value = "Status"

# Row 2 - "One"
$ws.Range("A2").Value = "FieldListOne"
$ws.Range("B2").Value = "OverrideListOne"
$ws.Range("C2").Value = "StatusListOne"

# Row 3 - "Two"
$ws.Range("A3").Value = "FieldListTwo"
$ws.Range("B3").Value = "OverrideListTwo"
$ws.Range("C3").Value = "StatusListTwo"

# Row 4 - "Three"
$ws.Range("A4").Value = "FieldListThree"
$ws.Range("B4").Value = "OverrideistThree"
$ws.Range("C4").Value = "StatusListThree"
